# Generate Report for Handback
#
# Refreshes the two "handback" file entries tracked by this report:
#   row 2: 869af6e9-4694-455d-9ad4-10fddd5aa9ba.md  -> 2097cada-03e0-4294-90f2-7cf176ed8740.md
#   row 3: c520e0be-cb63-4b87-b632-783e5fdae606.md  -> ffff0fb3cbc8-8b82-4f1d-8bcd-c84d07c6a43e.md
# along with the associated correspondence .xlf file names and the
# handoff/handback timestamps for each locale sheet.

$wb = $excel.ActiveWorkbook

$oldGuid1 = "869af6e9-4694-455d-9ad4-10fddd5aa9ba"
$newGuid1 = "2097cada-03e0-4294-90f2-7cf176ed8740"
$oldGuid2 = "c520e0be-cb63-4b87-b632-783e5fdae606"
$newGuid2 = "ffff0fb3cbc8-8b82-4f1d-8bcd-c84d07c6a43e"

$newXlfHash = "47ef2e3cb19969355ac3fa5826937b47bb3b32fa"

# ---------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A2").Value = "$newGuid1.md"
$wsOverview.Range("G2").Value = "2016-08-26 23:01:51"

$wsOverview.Range("A3").Value = "$newGuid2.md"
$wsOverview.Range("G3").Value = "2016-08-26 23:01:51"

# Hyperlink display text needs to track the renamed files too. This engine's
# Hyperlink objects can't be edited/deleted in place, so drop the sheet's
# hyperlinks and re-add them pointing at the same (unchanged) targets with
# the refreshed display text.
$wsOverview.Hyperlinks.Delete()
$wsOverview.Hyperlinks.Add(
    $wsOverview.Range("B2"),
    "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/5f66d4d15292c83e9c89be8d2fc7d5371ba2e48d/e2e/$oldGuid1.md",
    [System.Reflection.Missing]::Value,
    [System.Reflection.Missing]::Value,
    "e2e\$newGuid1.md"
)
$wsOverview.Hyperlinks.Add(
    $wsOverview.Range("B3"),
    "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/5f66d4d15292c83e9c89be8d2fc7d5371ba2e48d/e2e/$oldGuid2.md",
    [System.Reflection.Missing]::Value,
    [System.Reflection.Missing]::Value,
    "e2e\$newGuid2.md"
)

# ---------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$wsZhCn.Range("A2").Value = "$newGuid1.md"
$wsZhCn.Range("G2").Value = "$newGuid1.$newXlfHash.zh-cn.xlf"
$wsZhCn.Range("H2").Value = "2016-08-26 23:01:46"
$wsZhCn.Range("I2").Value = "$newGuid1.md"
$wsZhCn.Range("J2").Value = "$newGuid1.$newXlfHash.zh-cn.xlf"
$wsZhCn.Range("K2").Value = "2016-08-26 23:02:08"

$wsZhCn.Range("A3").Value = "$newGuid2.md"
$wsZhCn.Range("G3").Value = "$newGuid1.$newXlfHash.zh-cn.xlf"
$wsZhCn.Range("H3").Value = "2016-08-26 23:01:46"
$wsZhCn.Range("I3").Value = "$newGuid2.md"
$wsZhCn.Range("J3").Value = "$newGuid1.$newXlfHash.zh-cn.xlf"
$wsZhCn.Range("K3").Value = "2016-08-26 23:02:08"

$wsZhCn.Hyperlinks.Delete()
$wsZhCn.Hyperlinks.Add(
    $wsZhCn.Range("A2"),
    "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/5f66d4d15292c83e9c89be8d2fc7d5371ba2e48d/e2e/$oldGuid1.md",
    [System.Reflection.Missing]::Value,
    [System.Reflection.Missing]::Value,
    "$newGuid1.md"
)
$wsZhCn.Hyperlinks.Add(
    $wsZhCn.Range("I2"),
    "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/38729ae67f40fb308791b631bb102e6b292e690c/e2e/$oldGuid1.md",
    [System.Reflection.Missing]::Value,
    [System.Reflection.Missing]::Value,
    "$newGuid1.md"
)
$wsZhCn.Hyperlinks.Add(
    $wsZhCn.Range("A3"),
    "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/5f66d4d15292c83e9c89be8d2fc7d5371ba2e48d/e2e/$oldGuid2.md",
    [System.Reflection.Missing]::Value,
    [System.Reflection.Missing]::Value,
    "$newGuid2.md"
)
$wsZhCn.Hyperlinks.Add(
    $wsZhCn.Range("I3"),
    "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/38729ae67f40fb308791b631bb102e6b292e690c/e2e/$oldGuid2.md",
    [System.Reflection.Missing]::Value,
    [System.Reflection.Missing]::Value,
    "$newGuid2.md"
)

# ---------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsDeDe.Range("A2").Value = "$newGuid1.md"
$wsDeDe.Range("G2").Value = "$newGuid1.$newXlfHash.de-de.xlf"
$wsDeDe.Range("H2").Value = "2016-08-26 23:01:51"
$wsDeDe.Range("I2").Value = "$newGuid1.md"
$wsDeDe.Range("J2").Value = "$newGuid1.$newXlfHash.de-de.xlf"
$wsDeDe.Range("K2").Value = "2016-08-26 23:02:16"

$wsDeDe.Range("A3").Value = "$newGuid2.md"
$wsDeDe.Range("G3").Value = "$newGuid1.$newXlfHash.de-de.xlf"
$wsDeDe.Range("H3").Value = "2016-08-26 23:01:51"
$wsDeDe.Range("I3").Value = "$newGuid2.md"
$wsDeDe.Range("J3").Value = "$newGuid1.$newXlfHash.de-de.xlf"
$wsDeDe.Range("K3").Value = "2016-08-26 23:02:16"

$wsDeDe.Hyperlinks.Delete()
$wsDeDe.Hyperlinks.Add(
    $wsDeDe.Range("A2"),
    "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/5f66d4d15292c83e9c89be8d2fc7d5371ba2e48d/e2e/$oldGuid1.md",
    [System.Reflection.Missing]::Value,
    [System.Reflection.Missing]::Value,
    "$newGuid1.md"
)
$wsDeDe.Hyperlinks.Add(
    $wsDeDe.Range("I2"),
    "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/3bebbe89bef408967db0e94e5612c6013e5f1171/e2e/$oldGuid1.md",
    [System.Reflection.Missing]::Value,
    [System.Reflection.Missing]::Value,
    "$newGuid1.md"
)
$wsDeDe.Hyperlinks.Add(
    $wsDeDe.Range("A3"),
    "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/5f66d4d15292c83e9c89be8d2fc7d5371ba2e48d/e2e/$oldGuid2.md",
    [System.Reflection.Missing]::Value,
    [System.Reflection.Missing]::Value,
    "$newGuid2.md"
)
$wsDeDe.Hyperlinks.Add(
    $wsDeDe.Range("I3"),
    "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/3bebbe89bef408967db0e94e5612c6013e5f1171/e2e/$oldGuid2.md",
    [System.Reflection.Missing]::Value,
    [System.Reflection.Missing]::Value,
    "$newGuid2.md"
)
